$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Content.Find.Execute("2024-05-04 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-05 Sunday", 2) | Out-Null

# Update the arithmetic table cells (row-major order, matching document order)
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "4+9="  # was "84-18="
$tbl.Cell(1, 2).Range.Text = "17+25="  # was "26+6="
$tbl.Cell(1, 3).Range.Text = "90-67="  # was "39+4="
$tbl.Cell(1, 4).Range.Text = "69+23="  # was "77+8="
$tbl.Cell(1, 5).Range.Text = "48+33="  # was "75+17="
$tbl.Cell(2, 1).Range.Text = "75+6="  # was "27+7="
$tbl.Cell(2, 2).Range.Text = "49+28="  # was "86-68="
$tbl.Cell(2, 3).Range.Text = "19+16="  # was "29+13="
$tbl.Cell(2, 4).Range.Text = "72-45="  # was "80-5="
$tbl.Cell(2, 5).Range.Text = "74-69="  # was "29+4="
$tbl.Cell(3, 1).Range.Text = "44-5="  # was "60-37="
$tbl.Cell(3, 2).Range.Text = "23-18="  # was "90-18="
$tbl.Cell(3, 3).Range.Text = "88+6="  # was "65+17="
$tbl.Cell(3, 4).Range.Text = "16+66="  # was "34-8="
$tbl.Cell(3, 5).Range.Text = "93-78="  # was "94-9="
$tbl.Cell(4, 1).Range.Text = "93-24="  # was "36+7="
$tbl.Cell(4, 2).Range.Text = "81-62="  # was "66-49="
$tbl.Cell(4, 3).Range.Text = "52-8="  # was "42-7="
$tbl.Cell(4, 4).Range.Text = "80-43="  # was "8+54="
$tbl.Cell(4, 5).Range.Text = "37+38="  # was "92-89="
$tbl.Cell(5, 1).Range.Text = "12+49="  # was "43-8="
$tbl.Cell(5, 2).Range.Text = "15-9="  # was "87-69="
$tbl.Cell(5, 3).Range.Text = "76-29="  # was "39+44="
$tbl.Cell(5, 4).Range.Text = "94-89="  # was "86-77="
$tbl.Cell(6, 1).Range.Text = "77+17="  # was "36+28="
$tbl.Cell(6, 2).Range.Text = "86-19="  # was "57+28="
$tbl.Cell(6, 3).Range.Text = "84-5="  # was "17+77="
$tbl.Cell(6, 4).Range.Text = "68+3="  # was "91-74="
$tbl.Cell(6, 5).Range.Text = "71-65="  # was "71-15="
$tbl.Cell(7, 1).Range.Text = "62-14="  # was "80-13="
$tbl.Cell(7, 2).Range.Text = "83-44="  # was "15+16="
$tbl.Cell(7, 3).Range.Text = "93-78="  # was "75+16="
$tbl.Cell(7, 4).Range.Text = "52-14="  # was "24+49="
$tbl.Cell(7, 5).Range.Text = "8+6="  # was "52+9="
$tbl.Cell(8, 1).Range.Text = "68+8="  # was "19+5="
$tbl.Cell(8, 2).Range.Text = "48+36="  # was "17+39="
$tbl.Cell(8, 3).Range.Text = "62-3="  # was "55-38="
$tbl.Cell(8, 4).Range.Text = "26+59="  # was "20-11="
$tbl.Cell(8, 5).Range.Text = "62-13="  # was "39+17="
$tbl.Cell(9, 1).Range.Text = "38+48="  # was "78-9="
$tbl.Cell(9, 2).Range.Text = "25+7="  # was "86-39="
$tbl.Cell(9, 3).Range.Text = "59+35="  # was "78+7="
$tbl.Cell(9, 4).Range.Text = "27+65="  # was "73-34="
$tbl.Cell(9, 5).Range.Text = "47+9="  # was "36+25="
$tbl.Cell(10, 1).Range.Text = "71-25="  # was "74-58="
$tbl.Cell(10, 2).Range.Text = "83-68="  # was "53+38="
$tbl.Cell(10, 3).Range.Text = "71-25="  # was "15+27="
$tbl.Cell(10, 4).Range.Text = "92-56="  # was "32-6="
$tbl.Cell(10, 5).Range.Text = "48+23="  # was "59+24="
$tbl.Cell(11, 1).Range.Text = "55-29="  # was "90-51="
$tbl.Cell(11, 2).Range.Text = "46+46="  # was "7+85="
$tbl.Cell(11, 3).Range.Text = "39+5="  # was "45+7="
$tbl.Cell(11, 4).Range.Text = "49+42="  # was "74-26="
$tbl.Cell(11, 5).Range.Text = "28+35="  # was "93-8="
$tbl.Cell(12, 1).Range.Text = "16+46="  # was "72-6="
$tbl.Cell(12, 2).Range.Text = "71-29="  # was "81-22="
$tbl.Cell(12, 3).Range.Text = "44-17="  # was "19+65="
$tbl.Cell(12, 4).Range.Text = "17+66="  # was "81-33="
$tbl.Cell(12, 5).Range.Text = "82-19="  # was "39+26="
$tbl.Cell(13, 1).Range.Text = "38+26="  # was "6+56="
$tbl.Cell(13, 2).Range.Text = "56+5="  # was "77-69="
$tbl.Cell(13, 3).Range.Text = "13+49="  # was "85-49="
$tbl.Cell(13, 4).Range.Text = "9+19="  # was "27-18="
$tbl.Cell(13, 5).Range.Text = "91-85="  # was "9+88="
$tbl.Cell(14, 1).Range.Text = "56-8="  # was "43+8="
$tbl.Cell(14, 2).Range.Text = "60-33="  # was "31-18="
$tbl.Cell(14, 3).Range.Text = "37+37="  # was "70-28="
$tbl.Cell(14, 4).Range.Text = "43+29="  # was "15+67="
$tbl.Cell(14, 5).Range.Text = "19+37="  # was "94-19="
$tbl.Cell(15, 1).Range.Text = "69+29="  # was "82-59="
$tbl.Cell(15, 2).Range.Text = "9+76="  # was "91-16="
$tbl.Cell(15, 3).Range.Text = "26+6="  # was "39+57="
$tbl.Cell(15, 4).Range.Text = "53-34="  # was "14+29="
$tbl.Cell(15, 5).Range.Text = "27+35="  # was "92-17="
$tbl.Cell(16, 1).Range.Text = "88-9="  # was "80-49="
$tbl.Cell(16, 2).Range.Text = "93-34="  # was "15+26="
$tbl.Cell(16, 3).Range.Text = "8+45="  # was "82-29="
$tbl.Cell(16, 4).Range.Text = "63+9="  # was "31-14="
$tbl.Cell(16, 5).Range.Text = "29+14="  # was "67+14="
$tbl.Cell(17, 1).Range.Text = "61-7="  # was "8+18="
$tbl.Cell(17, 2).Range.Text = "37-18="  # was "94-59="
$tbl.Cell(17, 3).Range.Text = "25+67="  # was "27+39="
$tbl.Cell(17, 4).Range.Text = "91-33="  # was "25+48="
$tbl.Cell(17, 5).Range.Text = "30-7="  # was "80-57="
$tbl.Cell(18, 1).Range.Text = "18+46="  # was "29+43="
$tbl.Cell(18, 2).Range.Text = "48+15="  # was "94-66="
$tbl.Cell(18, 3).Range.Text = "25+16="  # was "78+5="
$tbl.Cell(18, 4).Range.Text = "71-43="  # was "25+36="
$tbl.Cell(18, 5).Range.Text = "44+7="  # was "94-59="
$tbl.Cell(19, 1).Range.Text = "39+46="  # was "47+45="
$tbl.Cell(19, 2).Range.Text = "82-23="  # was "17+59="
$tbl.Cell(19, 3).Range.Text = "63-37="  # was "96-9="
$tbl.Cell(19, 4).Range.Text = "72-66="  # was "39+42="
$tbl.Cell(19, 5).Range.Text = "90-82="  # was "39+42="
$tbl.Cell(20, 1).Range.Text = "9+9="  # was "64-48="
$tbl.Cell(20, 2).Range.Text = "41-15="  # was "59+5="
$tbl.Cell(20, 3).Range.Text = "7+25="  # was "27+69="
$tbl.Cell(20, 4).Range.Text = "95-56="  # was "9+65="
$tbl.Cell(20, 5).Range.Text = "85-76="  # was "44-15="

Write-Host "Done"
